$d = $word.ActiveDocument

# Locate the final (bookmarked) paragraph and the empty paragraph
# immediately above it, by walking back from the end of the document so
# the edit is resilient to unrelated content earlier in the body.
$finalPara = $d.Paragraphs.Last
$emptyParaB = $finalPara.Previous()

# Drop one of the two blank paragraphs that sit just above the paragraph
# that carries the "_GoBack" bookmark -- only one should remain.
$emptyParaB.Range.Delete()

# Re-fetch the (now immediately-preceding) bookmarked paragraph.
$finalPara = $d.Paragraphs.Last

# Add the trailing space *after* the bookmark first, while the paragraph
# is still empty -- inserting at the paragraph's end lands the new run
# after the bookmarkStart/bookmarkEnd pair already anchored there.
$tailRange = $d.Range($finalPara.Range.End, $finalPara.Range.End)
$tailRange.InsertAfter(" ")

# Now build up the explanatory sentence *before* the bookmark, one chunk
# at a time, inserting each new chunk immediately before the ones already
# placed at the paragraph's start. Working backwards (last sentence
# fragment first) keeps every chunk in its own run instead of being
# coalesced into a single run.
$finalPara = $d.Paragraphs.Last
$insertPos = $finalPara.Range.Start

$sentenceChunks = @(
    "pears in the receiver passband.",
    "wnlink ap",
    " to automatically change modes based on what type of do",
    "allow the radio",
    "input to the radio in order to ",
    "meaningful ",
    "Manual configuration is the baseline operation for Phase 4 Ground. This document describes an optional accessory shift-knob that supplies "
)

foreach ($chunk in $sentenceChunks) {
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertBefore($chunk)
}
